$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93, shifting existing rows 93:125 down to 94:126
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with its data
$ws.Cells.Item(93, 1).Value = 3
$ws.Cells.Item(93, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(93, 3).Value = "Coquimbo"
$ws.Cells.Item(93, 4).Value = 44463
$ws.Cells.Item(93, 5).Value = 5
$ws.Cells.Item(93, 6).Value = 100112010
$ws.Cells.Item(93, 7).Value = "Achicoria"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 115
$ws.Cells.Item(93, 11).Value = 5500
$ws.Cells.Item(93, 12).Value = 6000
$ws.Cells.Item(93, 13).Value = 5739
$ws.Cells.Item(93, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(93, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(93, 16).Value = 359
$ws.Cells.Item(93, 17).Value = 16
$ws.Cells.Item(93, 18).Value = "Hortaliza"

# Match the date cell style used by the rest of column D (s="2", datetime format)
$ws.Cells.Item(93, 4).NumberFormat = $ws.Cells.Item(94, 4).NumberFormat
